# Swap the species-record data between paired rows that were logged
# with swapped details (row 29<->30, 43<->44, 50<->51, 56<->57).
# Columns A,B,D,E,F,G,H,M,Q,R,Z,AB hold the per-record data that needs
# to be exchanged between each pair of rows; the rest of the row
# (location/observer metadata) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","M","Q","R","Z","AB")
$pairs = @(@(29,30), @(43,44), @(50,51), @(56,57))

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
